# Weekly update: a new price record for "Ají - Cacho cabra verde" is
# inserted at row 119 (Macroferia Regional de Talca, Hortaliza subset).
# This pushes the existing row 119 ("Cacho cabra verde", 2021-02-11) down
# to row 120, and the existing row 120 ("Americana (o)") down to row 121.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 119, shifting rows 119:120
# down to 120:121 (and extending the used range / dimension to R121).
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new weekly record.
$ws.Cells.Item(119, 1).Value = 5
$ws.Cells.Item(119, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(119, 3).Value = "Maule"
$ws.Cells.Item(119, 4).Value = 44448
$ws.Cells.Item(119, 5).Value = 7
$ws.Cells.Item(119, 6).Value = 100112021
$ws.Cells.Item(119, 7).Value = "Ají"
$ws.Cells.Item(119, 8).Value = "Cacho cabra verde"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 100
$ws.Cells.Item(119, 11).Value = 70000
$ws.Cells.Item(119, 12).Value = 70000
$ws.Cells.Item(119, 13).Value = 70000
$ws.Cells.Item(119, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(119, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(119, 16).Value = 2800
$ws.Cells.Item(119, 17).Value = 25
$ws.Cells.Item(119, 18).Value = "Hortaliza"
